# Adds a "Select Year / Current Estimate" sub-menu block under the
# "Estimates" menu entry, and an extra trailing blank paragraph after
# the "Admin" menu entry, per the source diff.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $target) {
    foreach ($p in $doc.Paragraphs) {
        $text = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($text -eq $target) {
            return $p
        }
    }
    return $null
}

function Insert-WordMLParagraphsAfter($doc, $paragraph, $paragraphsXml) {
    # Inserts one or more <w:p>...</w:p> fragments immediately after the
    # given paragraph by first materialising a placeholder paragraph
    # right after it, then replacing that placeholder's content via
    # Range.InsertXML so that tab characters round-trip as real
    # <w:tab/> run content (matching how Word itself persists them)
    # instead of literal "\t" inside <w:t>.
    if ($paragraph -eq $null) {
        throw "Insert-WordMLParagraphsAfter: target paragraph not found"
    }
    $idx = $paragraph.Index
    $r = $paragraph.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $placeholder = $doc.Paragraphs.Item($idx + 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $paragraphsXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $placeholder.Range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) New sub-menu paragraphs right after "Estimates", before "Stock Receive"
# ---------------------------------------------------------------------
$estimatesPara = Find-ParagraphByText $d "`tEstimates"

$newMenuParas = (
    '<w:p><w:r><w:tab/><w:t>Select Year</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Current Estimate</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> &gt; Change if not submitted &gt; Submit if not submitted</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/>' +
        '<w:t xml:space="preserve">When no previous estimate &gt; Create new Estimate for the year </w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/>' +
        '<w:t>Half Filled Estimate &gt; Load Previous Estimate</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>If Submitted &gt; Read-only Estimate</w:t></w:r></w:p>' +
    '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Search Estimates</w:t></w:r></w:p>'
)

Insert-WordMLParagraphsAfter $d $estimatesPara $newMenuParas

# InsertXML silently swallows a *trailing* empty <w:p/> when it is the
# last of several paragraphs in one fragment, so the blank separator
# paragraph before "Stock Receive" is added with its own call.
$searchEstimatesPara = Find-ParagraphByText $d "`t`tSearch Estimates"
Insert-WordMLParagraphsAfter $d $searchEstimatesPara '<w:p/>'

# ---------------------------------------------------------------------
# 2) One extra blank paragraph right after "Admin", before the existing
#    trailing blank paragraph / sectPr
# ---------------------------------------------------------------------
$adminPara = Find-ParagraphByText $d "`tAdmin"

Insert-WordMLParagraphsAfter $d $adminPara '<w:p/>'

Write-Host ("Final paragraph count: " + $d.Paragraphs.Count)
foreach ($pp in $d.Paragraphs) {
    Write-Host ("[" + $pp.Range.Text + "]")
}
